$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "may16_may_31" -> "may16_may31" (remove stray underscore before 31)
$ws.Range("A5").Value = "may16_may31"

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("C5:D6").Select()
